$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove row 2 ("H 72") entirely; subsequent rows shift up by one.
$ws.Rows("2:2").Delete()
